$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 335.7931
$ws.Range("I92").Value = 286.5926
$ws.Range("K92").Value = 286.5926
$ws.Range("M92").Value = 961.4074000000001
$ws.Range("H96").Value = 1512
$ws.Range("I96").Value = 760
$ws.Range("J96").Value = 3266.6667
$ws.Range("K96").Value = 2280
$ws.Range("L96").Value = 9800.000100000001
$ws.Range("M96").Value = -907
$ws.Range("N96").Value = -12546.0001
$ws.Range("H97").Value = 996.3333
$ws.Range("J97").Value = 994.5
$ws.Range("L97").Value = 2983.5
$ws.Range("N97").Value = -3975.5
$ws.Range("H100").Value = 3413.2
$ws.Range("I100").Value = 2867.0908
$ws.Range("J100").Value = 4915
$ws.Range("K100").Value = 2867.0908
$ws.Range("L100").Value = 4915
$ws.Range("M100").Value = -2326.0908
$ws.Range("N100").Value = -5997
$ws.Range("H101").Value = 1017
$ws.Range("I101").Value = 333.33334
$ws.Range("J101").Value = 2042.5
$ws.Range("K101").Value = 1000.00002
$ws.Range("L101").Value = 6127.5
$ws.Range("M101").Value = 621.9999799999999
$ws.Range("N101").Value = -9371.5
$ws.Range("H132").Value = 4652824
$ws.Range("I132").Value = 4763367.5
$ws.Range("K132").Value = 14290102.5
$ws.Range("M132").Value = -14287572.5
$ws.Range("H137").Value = 4764071
$ws.Range("I137").Value = 6668265.5
$ws.Range("J137").Value = 3583.8333
$ws.Range("K137").Value = 20004796.5
$ws.Range("L137").Value = 10751.4999
$ws.Range("M137").Value = -20002246.5
$ws.Range("N137").Value = -15851.4999

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3983.25
$ws.Range("I61").Value = 1701.8
$ws.Range("J61").Value = 5612.857
$ws.Range("K61").Value = 1701.8
$ws.Range("L61").Value = 5612.857
$ws.Range("M61").Value = -1489.8
$ws.Range("N61").Value = -6036.857
$ws.Range("H102").Value = 1893.238
$ws.Range("I102").Value = 1553.2222
$ws.Range("J102").Value = 3933.3333
$ws.Range("K102").Value = 1553.2222
$ws.Range("L102").Value = 3933.3333
$ws.Range("M102").Value = 68.77780000000007
$ws.Range("N102").Value = -7177.3333
$ws.Range("H132").Value = 18870746
$ws.Range("I132").Value = 21741590
$ws.Range("J132").Value = 5199.143
$ws.Range("K132").Value = 65224770
$ws.Range("L132").Value = 15597.429
$ws.Range("M132").Value = -65222240
$ws.Range("N132").Value = -20657.429
$ws.Range("H136").Value = 3983.25
$ws.Range("I136").Value = 1701.8
$ws.Range("J136").Value = 5612.857
$ws.Range("K136").Value = 5105.4
$ws.Range("L136").Value = 16838.571
$ws.Range("M136").Value = -2555.4
$ws.Range("N136").Value = -21938.571
$ws.Range("H138").Value = 50476.332
$ws.Range("J138").Value = 50476.332
$ws.Range("L138").Value = 50476.332
$ws.Range("N138").Value = -60756.332

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1707
$ws.Range("I86").Value = 1388
$ws.Range("J86").Value = 2309.5557
$ws.Range("K86").Value = 1388
$ws.Range("L86").Value = 2309.5557
$ws.Range("M86").Value = -265
$ws.Range("N86").Value = -4555.5557
$ws.Range("H89").Value = 1707
$ws.Range("I89").Value = 1388
$ws.Range("J89").Value = 2309.5557
$ws.Range("K89").Value = 6940
$ws.Range("L89").Value = 11547.7785
$ws.Range("M89").Value = -1324
$ws.Range("N89").Value = -22779.7785
$ws.Range("H94").Value = 470.27274
$ws.Range("I94").Value = 437.3
$ws.Range("J94").Value = 800
$ws.Range("K94").Value = 437.3
$ws.Range("L94").Value = 800
$ws.Range("M94").Value = 13.69999999999999
$ws.Range("N94").Value = -1702
$ws.Range("H107").Value = 4666.6665
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 4666.6665
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 4666.6665
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -8506.666499999999
$ws.Range("H133").Value = 19950
$ws.Range("J133").Value = 19950
$ws.Range("L133").Value = 19950
$ws.Range("N133").Value = -30070
$ws.Range("H134").Value = 4017.0334
$ws.Range("I134").Value = 3604.25
$ws.Range("J134").Value = 5668.1665
$ws.Range("K134").Value = 10812.75
$ws.Range("L134").Value = 17004.4995
$ws.Range("M134").Value = -8277.75
$ws.Range("N134").Value = -22074.4995

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H44").Value = 70071
$ws.Range("J44").Value = 70071
$ws.Range("L44").Value = 70071
$ws.Range("N44").Value = -70955
$ws.Range("H132").Value = 3495.0334
$ws.Range("I132").Value = 2089.4707
$ws.Range("K132").Value = 6268.4121
$ws.Range("M132").Value = -3738.4121
$ws.Range("H134").Value = 1935.8914
$ws.Range("I134").Value = 1612.1316
$ws.Range("J134").Value = 3473.75
$ws.Range("K134").Value = 4836.3948
$ws.Range("L134").Value = 10421.25
$ws.Range("M134").Value = -2301.3948
$ws.Range("N134").Value = -15491.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 7584.2856
$ws.Range("I3").Value = 272.5
$ws.Range("J3").Value = 17333.334
$ws.Range("K3").Value = 817.5
$ws.Range("L3").Value = 52000.00199999999
$ws.Range("M3").Value = -705.5
$ws.Range("N3").Value = -52224.00199999999
$ws.Range("H17").Value = 230.6923
$ws.Range("I17").Value = 239.08333
$ws.Range("J17").Value = 130
$ws.Range("K17").Value = 717.24999
$ws.Range("L17").Value = 390
$ws.Range("M17").Value = -548.24999
$ws.Range("N17").Value = -728

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5340.5884
$ws.Range("J122").Value = 4176.154
$ws.Range("L122").Value = 12528.462
$ws.Range("N122").Value = -17428.462
$ws.Range("H132").Value = 2640.5625
$ws.Range("I132").Value = 1786.9565
$ws.Range("K132").Value = 5360.8695
$ws.Range("M132").Value = -2830.8695

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2130.9
$ws.Range("I93").Value = 1941.5
$ws.Range("J93").Value = 2415
$ws.Range("K93").Value = 1941.5
$ws.Range("L93").Value = 2415
$ws.Range("M93").Value = -693.5
$ws.Range("N93").Value = -4911
$ws.Range("H132").Value = 2687.244
$ws.Range("I132").Value = 1715.5834
$ws.Range("K132").Value = 5146.7502
$ws.Range("M132").Value = -2616.7502
$ws.Range("H133").Value = 29250
$ws.Range("J133").Value = 29250
$ws.Range("L133").Value = 29250
$ws.Range("N133").Value = -34310

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1724.75
$ws.Range("I100").Value = 1724.75
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 3449.5
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -2908.5
$ws.Range("N100").ClearContents()
$ws.Range("H122").Value = 528406.4399999999
$ws.Range("I122").Value = 626688.9
$ws.Range("K122").Value = 1880066.7
$ws.Range("M122").Value = -1877616.7
$ws.Range("H132").Value = 392088.7
$ws.Range("I132").Value = 591500.5
$ws.Range("J132").Value = 15422
$ws.Range("K132").Value = 1774501.5
$ws.Range("L132").Value = 46266
$ws.Range("M132").Value = -1771971.5
$ws.Range("N132").Value = -51326
$ws.Range("H140").Value = 29616.125
$ws.Range("J140").Value = 29616.125
$ws.Range("L140").Value = 29616.125
$ws.Range("N140").Value = -39976.125
$ws.Range("H141").Value = 28333.334
$ws.Range("J141").Value = 28333.334
$ws.Range("L141").Value = 28333.334
$ws.Range("N141").Value = -38693.334
